# Updates the cryptos list (prices, volumes, and two row re-orderings)
# to match the data refresh performed by the GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C / E updates (plain text values; safe to assign directly) ---
$textUpdates = @{
    'E2'  = '  +0.87%  '
    'E3'  = '  +0.16%  '
    'E4'  = '  -0.13%  '
    'E5'  = '  +0.25%  '
    'E6'  = '  -0.09%  '
    'E7'  = '  -1.35%  '
    'E8'  = '  +0.05%  '
    'E9'  = '  -1.45%  '
    'E10' = '  +8.46%  '
    'E11' = '  -1.25%  '
    'E12' = '  +0.29%  '
    'E13' = '  -0.66%  '
    'E14' = '  -1.48%  '
    'E15' = '  -0.31%  '
    'E16' = '  +0.74%  '
    'E17' = '  -0.10%  '
    'E18' = '  +1.09%  '
    'E19' = '  -1.11%  '
    'E20' = '  -0.34%  '
    'E21' = '  -0.04%  '
    'E22' = '  -2.72%  '
    'E23' = '  -1.26%  '
    'B24' = 'Cosmos'
    'C24' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E24' = '  -0.50%  '
    'B25' = 'Monero'
    'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E25' = '  -0.40%  '
    'E26' = '  +1.17%  '
    'E27' = '  +1.38%  '
    'E28' = '  +2.16%  '
    'E29' = '  -0.68%  '
    'E30' = '  -3.05%  '
    'E31' = '  -1.62%  '
    'E32' = '  +0.37%  '
    'E33' = '  -0.15%  '
    'E34' = '  -2.31%  '
    'E35' = '  +0.47%  '
    'E36' = '  -1.75%  '
    'E37' = '  +2.07%  '
    'E38' = '  -3.29%  '
    'E39' = '  -1.19%  '
    'E40' = '  +0.07%  '
    'E41' = '  -0.08%  '
    'B42' = 'TheSandbox'
    'C42' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'E42' = '  -0.93%  '
    'B43' = 'FraxShare'
    'C43' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E43' = '  -0.29%  '
    'E44' = '  +2.16%  '
    'E45' = '  +2.57%  '
    'E46' = '  +0.97%  '
    'E47' = '  -1.85%  '
    'E48' = '  +1.69%  '
    'E49' = '  -1.04%  '
    'E50' = '  -2.18%  '
    'E51' = '  -0.48%  '
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# --- Column D (Price) updates ---
# These look like numbers (and some even look like valid floats), but the
# workbook stores them as literal text (t="inlineStr"/shared string), and a
# few of them (e.g. "9.250", "84.70", "30.567.36") would be silently mangled
# by Excel's automatic number conversion if assigned as plain strings.
# Forcing the cell to Text format before the assignment, then restoring the
# "Normal" style afterwards, preserves the exact text while leaving the
# cell's style index unchanged (so it stays unstyled, same as every other
# untouched Price cell).
$priceUpdates = @{
    'D2'  = '30.567.36'
    'D3'  = '1.863.32'
    'D4'  = '0.9984'
    'D5'  = '234.37'
    'D6'  = '0.9989'
    'D7'  = '0.4699'
    'D8'  = '0.2755'
    'D9'  = '0.06351'
    'D10' = '17.49'
    'D11' = '1.837.15'
    'D12' = '0.07455'
    'D13' = '4.968'
    'D14' = '84.70'
    'D15' = '0.6310'
    'D16' = '30.499.92'
    'D17' = '0.9990'
    'D18' = '234.63'
    'D19' = '12.68'
    'D20' = '0.000007351'
    'D21' = '0.9995'
    'D22' = '4.970'
    'D23' = '5.941'
    'D24' = '9.250'
    'D25' = '167.07'
    'D26' = '18.13'
    'D27' = '1.886'
    'D28' = '0.1030'
    'D29' = '1.375'
    'D30' = '4.101'
    'D31' = '3.857'
    'D32' = '0.04914'
    'D33' = '1.147'
    'D34' = '0.7081'
    'D35' = '2.704'
    'D36' = '0.01917'
    'D37' = '2.684'
    'D38' = '0.8795'
    'D39' = '1.966'
    'D40' = '105.80'
    'D41' = '0.9989'
    'D42' = '0.4085'
    'D43' = '5.541'
    'D44' = '7.222'
    'D45' = '0.1240'
    'D46' = '61.93'
    'D47' = '8.618'
    'D48' = '33.58'
    'D49' = '0.05555'
    'D50' = '1.374'
    'D51' = '0.3696'
}

foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}

Write-Output "Updated $($textUpdates.Count) label/volume cells and $($priceUpdates.Count) price cells"
